# Actualización automática 2025-06-09 11:15:09
#
# Updates the figures for client "PORCEKER S.A." (advisor CASTRO ALCIVAR EDA
# MARIA): an additional sale of 855.36 in the "PIEDRA SINTERIZADA" group is
# recorded, on top of the existing 1669.25 in "240X120 PORCELANATO"
# (bumped to 2166.91). This ripples through the monthly sales sheet (June
# column) and the monthly-compliance summary sheet, plus the "N de 54"
# tally label.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": per-client sales by product group
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 37 = PORCEKER S.A.
$wsGrupo.Range("C37").Value = 2166.91   # 240X120 PORCELANATO
$wsGrupo.Range("K37").Value = 855.36    # PIEDRA SINTERIZADA

# Tally row: one more client (now 3 of 54) sold PIEDRA SINTERIZADA
$wsGrupo.Range("K56").Value = "3 de 54"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": per-client sales by month
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 37 = PORCEKER S.A., column F = junio
$wsMensual.Range("F37").Value = 9376.389999999999
# Row 56 = grand total for junio
$wsMensual.Range("F56").Value = 24114.79

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": totals/compliance by product group
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column E (POR CUMPLIR) from 17 to 23 characters. Excel stores
# column widths with a +5/6 padding baked in, so back that out here to
# land exactly on a stored width of 23.
$wsCumplimiento.Range("E1").ColumnWidth = (23 - 5/6)

# Row 2 = 240X120 PORCELANATO
$wsCumplimiento.Range("D2").Value = 2166.91
$wsCumplimiento.Range("E2").Value = 3653.09
$wsCumplimiento.Range("F2").Value = 0.3723213058419244

# Row 15 = PIEDRA SINTERIZADA
$wsCumplimiento.Range("D15").Value = 8510.889999999999
$wsCumplimiento.Range("E15").Value = 7179.110000000001
$wsCumplimiento.Range("F15").Value = 0.542440407903123

# Row 19 = TOTAL
$wsCumplimiento.Range("D19").Value = 24210.55
$wsCumplimiento.Range("E19").Value = 66752.77900000001
$wsCumplimiento.Range("F19").Value = 0.2661572555243663
